$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain decimal numbers (e.g. "215.17",
# "0.0000252", "506.00"). The source data stores them as text, so force
# Text format on just those cells before assigning, which preserves exact
# digits/trailing zeros instead of Excel coercing them into numbers.
$textCells = @("D5", "D6", "D7", "D8", "D11", "D12", "D13", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '89.456.95'
$ws.Range("E2").Value = '  -2.85%  '
$ws.Range("D3").Value = '3.138.53'
$ws.Range("E3").Value = '  -3.52%  '
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").Value = '215.17'
$ws.Range("E5").Value = '  -1.36%  '
$ws.Range("D6").Value = '637.73'
$ws.Range("E6").Value = '  +2.18%  '
$ws.Range("D7").Value = '0.395'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = '0.772'
$ws.Range("E8").Value = '  +9.91%  '
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '3.134.81'
$ws.Range("E10").Value = '  -3.76%  '
$ws.Range("D11").Value = '0.559'
$ws.Range("E11").Value = '  -3.67%  '
$ws.Range("D12").Value = '0.179'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  -3.98%  '
$ws.Range("D14").Value = '5.34'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '89.202.34'
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("D16").Value = '3.707.66'
$ws.Range("E16").Value = '  -4.02%  '
$ws.Range("D17").Value = '32.35'
$ws.Range("E17").Value = '  -4.37%  '
$ws.Range("D18").Value = '3.118.13'
$ws.Range("E18").Value = '  -4.72%  '
$ws.Range("D19").Value = '3.40'
$ws.Range("E19").Value = '  +3.70%  '
$ws.Range("D20").Value = '0.0000228'
$ws.Range("E20").Value = '  +18.26%  '
$ws.Range("D21").Value = '13.27'
$ws.Range("E21").Value = '  -3.94%  '
$ws.Range("D22").Value = '426.82'
$ws.Range("E22").Value = '  -3.40%  '
$ws.Range("D23").Value = '8.37'
$ws.Range("E23").Value = '  -4.62%  '
$ws.Range("D24").Value = '4.92'
$ws.Range("E24").Value = '  -5.47%  '
$ws.Range("D25").Value = '5.50'
$ws.Range("E25").Value = '  +3.90%  '
$ws.Range("D26").Value = '82.07'
$ws.Range("E26").Value = '  +6.42%  '
$ws.Range("D27").Value = '11.59'
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").Value = '3.293.24'
$ws.Range("E28").Value = '  -4.57%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '0.157'
$ws.Range("E31").Value = '  -8.64%  '
$ws.Range("D32").Value = '4.04'
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").Value = '8.21'
$ws.Range("E33").Value = '  -5.27%  '
$ws.Range("D34").Value = '506.00'
$ws.Range("E34").Value = '  -7.62%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").Value = '7.02'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.143'
$ws.Range("E36").Value = '  +11.78%  '
$ws.Range("D37").Value = '1.29'
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("E38").Value = '  -4.47%  '
$ws.Range("D39").Value = '22.14'
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("D40").Value = '22.26'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("D44").Value = '0.365'
$ws.Range("E44").Value = '  -6.03%  '
$ws.Range("D45").Value = '145.47'
$ws.Range("E45").Value = '  -2.80%  '
$ws.Range("E46").Value = '  +3.18%  '
$ws.Range("D47").Value = '43.55'
$ws.Range("E47").Value = '  -3.86%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0661'
$ws.Range("E48").Value = '  +11.03%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '163.89'
$ws.Range("E49").Value = '  -7.57%  '
$ws.Range("D50").Value = '0.721'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").Value = '24.18'
$ws.Range("E51").Value = '  -2.44%  '
